# Updates the delivery-control sheet ("saida_controle"): refreshes the
# NºRomaneio/Cliente/Endereço/VOL data for the existing 8 orders (rows 2-9)
# and appends 6 new orders (rows 10-15), growing the used range from
# A1:G9 to A1:G15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("NºRomaneio") holds values that look numeric (e.g. "72.791") but
# must stay plain text, exactly like the rest of the sheet. Mark the range
# as Text *before* writing so Excel doesn't silently coerce them to numbers.
$ws.Range("C2:C15").NumberFormat = "@"

$ws.Cells.Item(2, 3).Value = '72.791'
$ws.Cells.Item(2, 4).Value = 'MASTER MEATS (BOUTIQUE DE CARNES)'
$ws.Cells.Item(2, 5).Value = 'Rua Abílio Soares, 731 - Paraíso São Paulo/SP CEP:04005003'
$ws.Cells.Item(2, 6).Value = 3

$ws.Cells.Item(3, 3).Value = '72.838'
$ws.Cells.Item(3, 4).Value = 'MYW1 O BAR LTDA (VASSOURA QUEBRADA - PERDIZES)'
$ws.Cells.Item(3, 5).Value = 'Rua Desembargador do Vale, 836, ANEXO 830 - Perdizes São Paulo/SP CEP:05010040'
$ws.Cells.Item(3, 6).Value = 10

$ws.Cells.Item(4, 3).Value = '72.859'
$ws.Cells.Item(4, 4).Value = 'MRL BUS COM. DE ALIM. EIRELI - ME (BUSGER - VILA MADALENA)'
$ws.Cells.Item(4, 5).Value = 'Rua Alves Guimarães, 1091, COZINHA 4 - Pinheiros São Paulo/SP CEP:05410-002'
$ws.Cells.Item(4, 6).Value = 15

$ws.Cells.Item(5, 3).Value = '72.891'
$ws.Cells.Item(5, 4).Value = 'ECULLY CHARBON RESTAURANTE LTDA (ECULLY CHARBON)'
$ws.Cells.Item(5, 5).Value = 'Rua Doutor Augusto de Miranda, 549 - Vila Pompéia São Paulo/SP CEP:05026000'
$ws.Cells.Item(5, 6).Value = 14

$ws.Cells.Item(6, 3).Value = '72.903'
$ws.Cells.Item(6, 4).Value = 'BAR & REST PICANHARIA DOS AMIGOS LTDA (PICANHARIA DOS AMIGOS - VILA LEOPOLDINA)'
$ws.Cells.Item(6, 5).Value = 'Rua Guaipá, 1017,  - Vila Leopoldina São Paulo/SP CEP:05089-001'
$ws.Cells.Item(6, 6).Value = 2

$ws.Cells.Item(7, 3).Value = '72.935'
$ws.Cells.Item(7, 4).Value = 'CAPITAO COM. E DIST. BEBIDAS E ALIMENTOS (CAPITAO BARLEY)'
$ws.Cells.Item(7, 5).Value = 'Rua Coriolano, 301 - Vila Romana São Paulo/SP CEP:05047001'
$ws.Cells.Item(7, 6).Value = 2

$ws.Cells.Item(8, 3).Value = '72.937'
$ws.Cells.Item(8, 4).Value = 'PARCEL SW BURGUER LTDA (N! BURGER - LAPA)'
$ws.Cells.Item(8, 5).Value = 'Rua Catão, 479, NBURGER - Vila Romana São Paulo/SP CEP:05049000'
$ws.Cells.Item(8, 6).Value = 11

$ws.Cells.Item(9, 3).Value = '72.941'
$ws.Cells.Item(9, 4).Value = 'TOSQUINHO LANCHES LTDA (TOSQUINHO LANCHES)'
$ws.Cells.Item(9, 5).Value = 'RUA CAMILO, 763, sem complemento - VILA ROMANA São Paulo/SP CEP:05045020'
$ws.Cells.Item(9, 6).Value = 8

$ws.Cells.Item(10, 3).Value = '72.956'
$ws.Cells.Item(10, 4).Value = '*CLIENTE AMOSTRA (CLIENTE AMOSTRA)*'
$ws.Cells.Item(10, 5).Value = 'Rua José Mariano Filho, 200,  - Jardim Oriental São Paulo/SP CEP:04347-180'
$ws.Cells.Item(10, 6).Value = 4

$ws.Cells.Item(11, 3).Value = '72.967'
$ws.Cells.Item(11, 4).Value = 'GILBERTO CAMPOS DE AZAMBUJA ME (ROYAL MEAT - PARAISO)'
$ws.Cells.Item(11, 5).Value = 'Rua Doutor Tomás Carvalhal, 626 - Paraíso São Paulo/SP CEP:04006001'
$ws.Cells.Item(11, 6).Value = 5

$ws.Cells.Item(12, 3).Value = '72.970'
$ws.Cells.Item(12, 4).Value = 'THE BEAR BURGER REST. LTDA EPP (THE BEAR BURGER)'
$ws.Cells.Item(12, 5).Value = 'Rua Caraíbas, 964, IMOBILIARIA ESTEVAM - Perdizes São Paulo/SP CEP:05020000'
$ws.Cells.Item(12, 6).Value = 5

$ws.Cells.Item(13, 3).Value = '72.988'
$ws.Cells.Item(13, 4).Value = 'BUSGER COM. DE ALIM. LTDA (BUSGER - KLABIN)'
$ws.Cells.Item(13, 5).Value = 'Rua Vergueiro, 4289,  - Vila Mariana São Paulo/SP CEP:04101-901'
$ws.Cells.Item(13, 6).Value = 30

$ws.Cells.Item(14, 3).Value = '72.990'
$ws.Cells.Item(14, 4).Value = 'ESTEFOODS COM. DE ALIMENTOS LTDA (BUSGER - BORGES LAGOA)'
$ws.Cells.Item(14, 5).Value = 'Rua Borges Lagoa, 1050,  - Vila Clementino São Paulo/SP CEP:04038-002'
$ws.Cells.Item(14, 6).Value = 26

$ws.Cells.Item(15, 3).Value = '73.008'
$ws.Cells.Item(15, 4).Value = 'BORGER BURGER LTDA (BORGER - PERDIZES)'
$ws.Cells.Item(15, 5).Value = 'Rua Cardoso de Almeida, 587,  - Perdizes São Paulo/SP CEP:05013-000'
$ws.Cells.Item(15, 6).Value = 2
